$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make D8's format match the rest of column D (date format) so the cell
# exists with the right style once its value is cleared.
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats

# Clear contents of A2:D8 but keep formatting
$ws.Range("A2:D8").ClearContents()

# Update D1 text (was "row4" via shared string index 9, now referencing the deduped index)
$ws.Range("D1").Value = "row4"

# Update selection to E8
$ws.Range("E8").Select()
